$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- Paragraph 1: merge "Internationalized Resource " + "Identifier" runs into one run ---
$para1 = $tr.Paragraphs(1, 1)
$merge1 = $para1.Characters(67, 37)
$merge1.Text = "Internationalized Resource Identifier"

# --- Paragraph 2: merge "Resource Description " + "Framework" runs into one run ---
$para2 = $tr.Paragraphs(2, 1)
$merge2 = $para2.Characters(6, 30)
$merge2.Text = "Resource Description Framework"

# --- Paragraph 3: reshape "Web Ontology " + "Language)" into "Web Ontology Language" + ") - " ---
# and append the new Russian sentence describing the purpose of OWL.
$para3 = $tr.Paragraphs(3, 1)
$tail = $para3.Characters(38, 22)
$tail.Text = "Web Ontology Language) – используются для описания структуры предметной области в семантической паутине."

# Re-assert italics on ") - " to split it from "Web Ontology Language" as its own run
$dashPart = $para3.Characters(59, 4)
$dashPart.Font.Italic = 1

# Turn off italics for the new Russian sentence, in two runs as per the source deck
$ruPart1 = $para3.Characters(63, 57)
$ruPart1.Font.Italic = 0

$ruPart2 = $para3.Characters(120, 22)
$ruPart2.Font.Italic = 0
